# Edit LOM3004.xlsx worksheet to match the target revision:
#  - Insert a new row at position 13 (pushing existing rows 13-23 down to 14-24)
#  - Update several cell contents with the revised Portuguese course text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row before the current row 13 ("Programa resumido:").
#    This shifts rows 13-23 down to 14-24, matching the new dimension A1:C24.
$ws.Rows("13:13").Insert()

# 2. Objetivos (row 10) - replace the placeholder professor name with the
#    actual course objectives text (PT).
$objetivosPt = "1. Descrever os principais processos de conformação mecânica utilizados na indústria metal mecânica.2. Munir o aluno de conhecimentos suficientes para especificar equipamentos e acessórios, usados no processo de conformação, com base nas solicitações mecânicas e variáveis do processo. 3. Ensinar ao aluno a definir tecnicamente o processo adequado de conformação mecânica de produtos da indústria metal mecânica."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# 3. Docentes responsáveis (new row 13) - fill in the professor's name that
#    used to sit under "Programa resumido:". The freshly inserted row has no
#    cells yet, so A13 must stay empty (no label on this line) and B13/C13
#    need to pick up the normal (non-bold, wrapped) value styles used
#    elsewhere in columns B/C rather than the bold style the insert applied.
$docente = "5840793 - Sérgio Schneider"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
$ws.Range("A13").Clear()
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true

# 4. Programa resumido (row 14) - replace "Semestral" with the short syllabus text (PT).
$resumidoPt = "1. Introdução à Teoria de Plasticidade.2. Classificação dos Processos de Conformação Mecânica.3. Metalurgia da Conformação Mecânica.4. Mecânica da Conformação: Forjamento, Trefilação, Extrusão, Laminação, Dobramento e Estampagem.5. Descrição dos Processos de Conformação Plástica6. Processamento de Chapas Metálicas: Estampagem, Dobramento, Calandragem e corte."
$ws.Range("B14").Value = $resumidoPt
$ws.Range("C14").Value = $resumidoPt

# 5. Programa (row 16) - replace the stray date with the full syllabus text (PT).
$programaPt = "1) Introdução à Teoria de Plasticidade: Análise de tensão e deformação 3D, tensões principais, critérios de escoamento, relação tensão-deformação no regime plástico, tensão efetiva e deformação efetiva, energia de deformação e trabalho plástico. 2) Classificação dos processos de conformação. Forjamento, Trefilação, Extrusão Laminação, Estiramento, Estampagem e Dobramento. 3) Metalurgia da Conformação: Noções básicas sobre recozimento, encruamento, conformabilidade, textura e anisotropia. Trabalho a frio e a quente.4) Mecânica da Conformação: métodos de cálculo, efeitos do atrito na conformação e noções básicas de lubrificação. Efeito da taxa de deformação. 5) Descrição dos Processos de Conformação Plástica. Forjamento: Trefilação, Extrusão e Laminação. Equipamentos e acessórios: características e noções de projeto e dimensionamento dos mesmos. Cálculo de carga desses processos.  Análise de defeitos que podem ocorrer nesses processos. Relações geométricas na laminação. Potência e torque de laminação.  6) Processamento de Chapas Metálica. Estampagem: Ensaios de Estampabilidade, Curva Limite de Conformação (CLC). Dobramento: tipos de dobramento, efeito mola, equipamentos e matrizes. Calandragem: tipos de calandragem e equipamentos. Corte: equipamentos de corte e aplicações."
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# 6. Método (row 19) - now holds the text that used to be under "Critério:".
$metodoTxt = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2)."
$ws.Range("B19").Value = $metodoTxt
$ws.Range("C19").Value = $metodoTxt

# 7. Critério (row 20) - now holds the text that used to be under "Norma de recuperação:".
$criterioTxt = "NS = (P1 + P2)/2Serão considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("B20").Value = $criterioTxt
$ws.Range("C20").Value = $criterioTxt

# 8. Norma de recuperação (row 21) - now holds the text that used to be under "Bibliografia:".
$normaTxt = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma:NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("B21").Value = $normaTxt
$ws.Range("C21").Value = $normaTxt

# 9. Bibliografia (row 22) - replace with the real bibliography list (PT).
$biblioTxt = ".  LARKE, E.C. The Rolling of Strip, Sheet, and Plate, Chapman and Hall, 19672.  HONEYCOMBE, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1968.3.  HOSFORD, W.F. Metal Forming: Mechanics and Metallurgy, Prentice-Hall, 1983. 4.  WEERTMAN, J. Elementary Dislocation Theory, Collier-McMillan, 1965. 5.  AVITZUR, B. Metal Forming: Processes and Analysis, McGraw-Hill, 1968.6. BRESCIANI Filho, E. e outros. Conformação Plástica dos Metais, Editora da UNICAMP           Campinas, Volumes 1 e 2, 1986. 7. CETLIN, P. R.; HELMAN, H. Fundamentos de Conformação Mecânica dos Metais. Art Liber: São Paulo, 2005.8. ROWE, G.W. Elements of Metalworking Theory. Edward Arnold Publishers, 1979.9.  JOHNSON, W.; MELLOR, P.B. Engineering Plasticity, Van Nostrand Reinhold, 1973. 10. DIETER, G. E. Metalurgia Mecânica. Guanabara Dois, 1981.11. SCHAEFFER, L. Introdução à Conformação Mecânica dos Metais, Ed. da UFRGS, 1983. 12. RODRIGUES, J. Tecnologia Mecânica. Volumes 1 e 2, Ed. Escolar, 2005. 13. CALLISTER, W. D. Ciência e Engenharia dos Materiais: Uma Introdução. Rio de Janeiro: LTC, 1999."
$ws.Range("B22").Value = $biblioTxt
$ws.Range("C22").Value = $biblioTxt
